$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before G (so old G..J shift to H..K)
$ws.Range("G1").EntireColumn.Insert()

# New header for column G
$ws.Range("G1").Value = "Alamat"

# New address values per row (mapped by assistant name)
$ws.Range("G4").Value = "Jl Srengseng Raya 88, Dki Jakarta"
$ws.Range("G2").Value = "Jl Prof Dr Sudarto 126 A, Jawa Tengah"
$ws.Range("G3").Value = "JL Pemuda No.1-G, Rawamangun"
